# Regional Availability Factor.xlsx - "updated 4.0 files and mdl"
#
# Intentional data/view edits captured by this commit:
#   1. About!C1  - "last updated" date bumped 2024-03-15 -> 2024-03-28 (serial 45366 -> 45379)
#   2. RAF-capacity!B24 (hydrogen combustion turbine capacity credit multiplier): 0.3 -> 1
#   3. RAF-capacity!B25 (hydrogen combined cycle capacity credit multiplier):     0.3 -> 1
#   4. The workbook's active/selected tab moves from "RAF-generation" to "RAF-capacity"
#   5. On RAF-capacity: selection moves to B25, the view is scrolled down (top row ~14)
#      and zoomed to 80%, and column A is narrowed a bit to fit its contents.

$wb = $excel.ActiveWorkbook

# --- 1. About sheet: bump the "last updated" date -------------------------
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45379

# --- 2 & 3. RAF-capacity sheet: raise hydrogen plant capacity credit ------
$wsCap = $wb.Worksheets.Item("RAF-capacity")
$wsCap.Range("B24").Value = 1
$wsCap.Range("B25").Value = 1

# --- 4. Make RAF-capacity the active sheet (was RAF-generation) ----------
$wsCap.Activate()

# --- 5. Update the RAF-capacity sheet view: column width, scroll position,
#        zoom level and selected cell -------------------------------------
$wsCap.Columns.Item(1).ColumnWidth = 28.1

[void]$wsCap.Range("B25").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 14
$win.ScrollColumn = 1
$win.Zoom = 80
